$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sessions")
$ws.Columns("L").Insert()
Write-Output "done"
